# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled run).
# Most D/E cells only get their text re-stamped with the latest scraped
# price / 1h-volume-change strings. A handful of D-column prices look like
# plain decimals (e.g. "1.000", "0.4295") which Excel would otherwise
# auto-coerce to numbers, so those cells are pre-formatted as Text ("@")
# before the value is written, to keep them as literal strings like the
# rest of the sheet (which stores everything as text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.299.50"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.872.40"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7119"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.70"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3109"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07701"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08361"
$ws.Range("D12").Value = "1.872.62"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("E14").Value = "  -1.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.31"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "29.312.28"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008264"
$ws.Range("E17").Value = "  +5.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.937"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.85"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("D20").Value = "2.132.81"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.848"
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1625"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.16"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.008"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.52"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.503"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("E31").Value = "  +5.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.284"
$ws.Range("E32").Value = "  -4.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05247"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.927"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7532"
$ws.Range("E35").Value = "  +3.27%  "
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.681"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01858"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.715"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").Value = "1.153.75"
$ws.Range("E40").Value = "  -1.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.361"
$ws.Range("E41").Value = "  +4.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.01"
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8879"
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.62"
$ws.Range("E44").Value = "  +2.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9997"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "2.029.04"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5192"
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.379"
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000120"
$ws.Range("E50").Value = "  +2.69%  "
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4295"
$ws.Range("E51").Value = "  +0.56%  "
